# Updated cryptos list values (Price column D, Volume(1h) column E).
#
# The source cells are plain text (t="inlineStr" in the OOXML) even though
# many Price values look numeric (e.g. "4.330", "207.08"). Writing such a
# string straight to Range.Value lets Excel's COM layer auto-detect it as a
# number (e.g. "4.330" -> 4.33, silently dropping the trailing zero and
# changing the cell's stored type). To keep these values as literal text we
# write them with a leading apostrophe (Excel's standard "force text" quote
# prefix) and then reset the cell's Style back to "Normal" so the
# quote-prefix formatting introduced by that apostrophe isn't left behind -
# this matches the original cells, which carry no explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '25.727.23' },
    @{ Cell = 'E2'; Value = '  -2.03%  ' },
    @{ Cell = 'D3'; Value = '1.608.84' },
    @{ Cell = 'E3'; Value = '  -4.06%  ' },
    @{ Cell = 'E4'; Value = '  +0.04%  ' },
    @{ Cell = 'D5'; Value = '207.08' },
    @{ Cell = 'E5'; Value = '  -2.56%  ' },
    @{ Cell = 'D6'; Value = '0.5161' },
    @{ Cell = 'E6'; Value = '  -2.20%  ' },
    @{ Cell = 'D7'; Value = '1.003' },
    @{ Cell = 'E7'; Value = '  +0.02%  ' },
    @{ Cell = 'D8'; Value = '0.2547' },
    @{ Cell = 'E8'; Value = '  -4.21%  ' },
    @{ Cell = 'D9'; Value = '0.06167' },
    @{ Cell = 'E9'; Value = '  -2.04%  ' },
    @{ Cell = 'E10'; Value = '  -6.31%  ' },
    @{ Cell = 'D11'; Value = '0.07514' },
    @{ Cell = 'E11'; Value = '  -0.59%  ' },
    @{ Cell = 'D12'; Value = '1.597.71' },
    @{ Cell = 'E12'; Value = '  -7.27%  ' },
    @{ Cell = 'D13'; Value = '4.330' },
    @{ Cell = 'E13'; Value = '  -3.21%  ' },
    @{ Cell = 'D14'; Value = '1.831.66' },
    @{ Cell = 'E14'; Value = '  -4.08%  ' },
    @{ Cell = 'D15'; Value = '0.5392' },
    @{ Cell = 'E15'; Value = '  -4.47%  ' },
    @{ Cell = 'D16'; Value = '0.0₅7765' },
    @{ Cell = 'E16'; Value = '  -3.53%  ' },
    @{ Cell = 'D17'; Value = '63.35' },
    @{ Cell = 'E17'; Value = '  -5.69%  ' },
    @{ Cell = 'D18'; Value = '25.720.10' },
    @{ Cell = 'E18'; Value = '  -1.32%  ' },
    @{ Cell = 'D19'; Value = '1.003' },
    @{ Cell = 'E19'; Value = '  +0.03%  ' },
    @{ Cell = 'D20'; Value = '4.577' },
    @{ Cell = 'E20'; Value = '  -5.40%  ' },
    @{ Cell = 'D21'; Value = '182.15' },
    @{ Cell = 'E21'; Value = '  -3.25%  ' },
    @{ Cell = 'E22'; Value = '  -4.76%  ' },
    @{ Cell = 'E23'; Value = '  +0.09%  ' },
    @{ Cell = 'D24'; Value = '5.991' },
    @{ Cell = 'E24'; Value = '  -3.74%  ' },
    @{ Cell = 'D25'; Value = '143.88' },
    @{ Cell = 'E25'; Value = '  -4.14%  ' },
    @{ Cell = 'D26'; Value = '0.1195' },
    @{ Cell = 'E26'; Value = '  -4.95%  ' },
    @{ Cell = 'D27'; Value = '7.295' },
    @{ Cell = 'E27'; Value = '  -4.14%  ' },
    @{ Cell = 'D28'; Value = '15.35' },
    @{ Cell = 'E28'; Value = '  -4.04%  ' },
    @{ Cell = 'D29'; Value = '1.352' },
    @{ Cell = 'E29'; Value = '  -0.91%  ' },
    @{ Cell = 'D30'; Value = '0.05820' },
    @{ Cell = 'E30'; Value = '  -6.68%  ' },
    @{ Cell = 'D31'; Value = '1.231' },
    @{ Cell = 'E31'; Value = '  -4.33%  ' },
    @{ Cell = 'D32'; Value = '3.333' },
    @{ Cell = 'E32'; Value = '  -5.41%  ' },
    @{ Cell = 'D33'; Value = '3.295' },
    @{ Cell = 'E33'; Value = '  -4.29%  ' },
    @{ Cell = 'D34'; Value = '1.581' },
    @{ Cell = 'E34'; Value = '  -3.69%  ' },
    @{ Cell = 'D35'; Value = '0.9568' },
    @{ Cell = 'E35'; Value = '  -4.81%  ' },
    @{ Cell = 'E36'; Value = '  -1.15%  ' },
    @{ Cell = 'D37'; Value = '2.698' },
    @{ Cell = 'E37'; Value = '  -1.44%  ' },
    @{ Cell = 'D38'; Value = '0.5696' },
    @{ Cell = 'E38'; Value = '  -6.25%  ' },
    @{ Cell = 'D39'; Value = '0.01572' },
    @{ Cell = 'E39'; Value = '  -3.20%  ' },
    @{ Cell = 'D40'; Value = '1.002' },
    @{ Cell = 'E40'; Value = '  -0.47%  ' },
    @{ Cell = 'D41'; Value = '0.8334' },
    @{ Cell = 'E41'; Value = '  -4.45%  ' },
    @{ Cell = 'D42'; Value = '5.611' },
    @{ Cell = 'E42'; Value = '  -8.27%  ' },
    @{ Cell = 'D43'; Value = '1.013.57' },
    @{ Cell = 'E43'; Value = '  -8.19%  ' },
    @{ Cell = 'D44'; Value = '98.64' },
    @{ Cell = 'E44'; Value = '  -1.39%  ' },
    @{ Cell = 'D45'; Value = '1.758.67' },
    @{ Cell = 'E45'; Value = '  -3.76%  ' },
    @{ Cell = 'E46'; Value = '  -2.92%  ' },
    @{ Cell = 'D47'; Value = '1.002' },
    @{ Cell = 'E47'; Value = '  -0.47%  ' },
    @{ Cell = 'D48'; Value = '53.60' },
    @{ Cell = 'E48'; Value = '  -4.82%  ' },
    @{ Cell = 'E49'; Value = '  -1.55%  ' },
    @{ Cell = 'D50'; Value = '7.797' },
    @{ Cell = 'E50'; Value = '  -2.81%  ' },
    @{ Cell = 'D51'; Value = '0.4210' },
    @{ Cell = 'E51'; Value = '  -1.08%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "'" + $u.Value
    $range.Style = "Normal"
}
